$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Header row: rename C1 and add a new D1 header ("QQ号（必填）").
#    Copy the format from the existing header cell (C1) onto D1 first, so the
#    new cell picks up the same style (bold/empty font xf) as its neighbours,
#    then change the text of C1 and set the text of D1.
# ---------------------------------------------------------------------------
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C1").Value = "序列号（必填）"
$ws.Range("D1").Value = "QQ号（必填）"

# ---------------------------------------------------------------------------
# 2. Row 20 ("1874" submitter): mark the row as struck-through with a grey
#    fill, and clear out the submitted serial number in C20.
#    Build the two needed combinations (with / without the date number
#    format) on scratch cells first, then paste the *formats only* onto the
#    real cells - this re-uses a single pair of style entries instead of
#    generating throw-away intermediate styles.
# ---------------------------------------------------------------------------
$scratchDate = $ws.Range("F20")
$scratchDate.NumberFormat = "yyyy/m/d h:mm:ss;@"
$scratchDate.Font.Strikethrough = $true
$scratchDate.Interior.Color = 14211288

$scratchPlain = $ws.Range("G20")
$scratchPlain.Font.Strikethrough = $true
$scratchPlain.Interior.Color = 14211288

$scratchDate.Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null

$scratchPlain.Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("F20:G20").Clear() | Out-Null

$ws.Range("C20").Value = ""

# ---------------------------------------------------------------------------
# 3. Append new submission rows (25-27).
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "李智杰"
$ws.Range("B25").Value = 45897.9370717593
$ws.Range("C25").Value = "c7d4d17b"

$ws.Range("A26").Value = [string][char]0x3000
$ws.Range("B26").Value = 45899.356400463
$ws.Range("C26").Value = "737c6559"

$ws.Range("A27").Value = "Golden"
$ws.Range("B27").Value = 45899.747337963
$ws.Range("C27").Value = "87227784"

# The new "QQ号" column (D) holds numeric-looking strings that must stay
# text. Force text formatting while entering them, then restore the default
# (un-styled) look by pasting the formatting of an existing plain cell.
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1535752313"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1123084248"

$ws.Range("A2").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
